$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dists")
$ws.Range("D17").Value = 14.6438946363933
Write-Output "done"
